$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-6,
# incrementing the serial date value from 45207 to 45208 (2023-10-08 -> 2023-10-09).
$ws.Range("C2:C6").Value = 45208
